# Replace the 15 lattice-multiplication exercise cells (5 rows x 3 cols)
# with new problems/partial-product scaffolding, per the target revision.
#
# Each cell's text is: "A x B" <br> "  d1    d2" <br> "  ----" <br>
#                       "d|    |" <br> "d|    |"
# where <br> is Word's manual-line-break character (vertical tab, chr(11)).
# Writing the whole cell Range.Text in one shot preserves the existing
# run formatting (sz=32) and regenerates the <w:br/> run breaks correctly.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)
$nl = [char]11

$t.Cell(1,1).Range.Text = "13 x 75" + $nl + "  7    5" + $nl + "  ----" + $nl + "1|    |" + $nl + "3|    |"
$t.Cell(1,2).Range.Text = "70 x 89" + $nl + "  8    9" + $nl + "  ----" + $nl + "7|    |" + $nl + "0|    |"
$t.Cell(1,3).Range.Text = "72 x 63" + $nl + "  6    3" + $nl + "  ----" + $nl + "7|    |" + $nl + "2|    |"

$t.Cell(2,1).Range.Text = "64 x 34" + $nl + "  3    4" + $nl + "  ----" + $nl + "6|    |" + $nl + "4|    |"
$t.Cell(2,2).Range.Text = "83 x 98" + $nl + "  9    8" + $nl + "  ----" + $nl + "8|    |" + $nl + "3|    |"
$t.Cell(2,3).Range.Text = "36 x 44" + $nl + "  4    4" + $nl + "  ----" + $nl + "3|    |" + $nl + "6|    |"

$t.Cell(3,1).Range.Text = "42 x 45" + $nl + "  4    5" + $nl + "  ----" + $nl + "4|    |" + $nl + "2|    |"
$t.Cell(3,2).Range.Text = "54 x 64" + $nl + "  6    4" + $nl + "  ----" + $nl + "5|    |" + $nl + "4|    |"
$t.Cell(3,3).Range.Text = "21 x 59" + $nl + "  5    9" + $nl + "  ----" + $nl + "2|    |" + $nl + "1|    |"

$t.Cell(4,1).Range.Text = "23 x 63" + $nl + "  6    3" + $nl + "  ----" + $nl + "2|    |" + $nl + "3|    |"
$t.Cell(4,2).Range.Text = "36 x 46" + $nl + "  4    6" + $nl + "  ----" + $nl + "3|    |" + $nl + "6|    |"
$t.Cell(4,3).Range.Text = "16 x 40" + $nl + "  4    0" + $nl + "  ----" + $nl + "1|    |" + $nl + "6|    |"

$t.Cell(5,1).Range.Text = "38 x 88" + $nl + "  8    8" + $nl + "  ----" + $nl + "3|    |" + $nl + "8|    |"
$t.Cell(5,2).Range.Text = "84 x 29" + $nl + "  2    9" + $nl + "  ----" + $nl + "8|    |" + $nl + "4|    |"
# Final cell in the target ends with a trailing break and no closing
# second partial-product line (matches the source revision exactly).
$t.Cell(5,3).Range.Text = "37 x 17" + $nl + "  1    7" + $nl + "  ----" + $nl + "3|    |" + $nl
